$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The refreshed Power Query dropped the "Progress" column's values (B2:B15).
$ws.Range("B2:B15").ClearContents()

# 2. The refreshed query table re-applies formatting to the "Trial Name" data cells.
$ws.Range("A2:A16").NumberFormat = "General"

# 3. A new 16th row (blank trial name) came back from the refreshed query.
$ws.Cells.Item(16, 1).Value = ""

# 4. Resize the Query1 table to the new extent A1:B16.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B16"))

# 5. Update the ExternalData_1 defined name to the new range.
$nm = $wb.Names.Item("ExternalData_1")
$nm.RefersTo = "=Sheet1!`$A`$1:`$B`$16"

Write-Host "done"
